$wb = $excel.ActiveWorkbook

function Set-CellValue {
    param($ws, $row, $col, $value)
    $ws.Cells.Item($row, $col).Value = $value
}

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 19 8 276.4516
Set-CellValue $ws 19 9 228.66667
Set-CellValue $ws 19 10 321.25
Set-CellValue $ws 19 11 228.66667
Set-CellValue $ws 19 12 321.25
Set-CellValue $ws 19 13 -53.66667000000001
Set-CellValue $ws 19 14 -671.25

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 98 8 3459.1
Set-CellValue $ws 98 9 3974.125
Set-CellValue $ws 98 11 3974.125
Set-CellValue $ws 98 13 -2476.125

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 113 8 899.5
Set-CellValue $ws 113 9 899.5
Set-CellValue $ws 113 11 899.5
Set-CellValue $ws 113 13 2354.5

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 122 8 3459.1
Set-CellValue $ws 122 9 3974.125
Set-CellValue $ws 122 11 11922.375
Set-CellValue $ws 122 13 -9472.375

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
Set-CellValue $ws 137 8 1359.6904
Set-CellValue $ws 137 9 1085.0358
Set-CellValue $ws 137 10 1909
Set-CellValue $ws 137 11 3255.1074
Set-CellValue $ws 137 12 5727
Set-CellValue $ws 137 13 -705.1074000000003
Set-CellValue $ws 137 14 -10827

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 2 8 647.7742
Set-CellValue $ws 2 9 430.03845
Set-CellValue $ws 2 10 1780
Set-CellValue $ws 2 11 430.03845
Set-CellValue $ws 2 12 1780
Set-CellValue $ws 2 13 -317.03845
Set-CellValue $ws 2 14 -2006

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 32 8 6820.3726
Set-CellValue $ws 32 9 5595
Set-CellValue $ws 32 10 12370.588
Set-CellValue $ws 32 11 5595
Set-CellValue $ws 32 12 12370.588
Set-CellValue $ws 32 13 -5308
Set-CellValue $ws 32 14 -12944.588

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 116 8 647.7742
Set-CellValue $ws 116 9 430.03845
Set-CellValue $ws 116 10 1780
Set-CellValue $ws 116 11 430.03845
Set-CellValue $ws 116 12 1780
Set-CellValue $ws 116 13 1863.96155
Set-CellValue $ws 116 14 -6368

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 122 8 2020.9445
Set-CellValue $ws 122 9 1629.9333
Set-CellValue $ws 122 10 3976
Set-CellValue $ws 122 11 4889.7999
Set-CellValue $ws 122 12 11928
Set-CellValue $ws 122 13 -2439.7999
Set-CellValue $ws 122 14 -16828

# Sheet ARM, row 133
$ws = $wb.Worksheets.Item("ARM")
Set-CellValue $ws 133 8 34936
Set-CellValue $ws 133 10 34936
Set-CellValue $ws 133 12 34936
Set-CellValue $ws 133 14 -39996

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 3 8 647.7742
Set-CellValue $ws 3 9 430.03845
Set-CellValue $ws 3 10 1780
Set-CellValue $ws 3 11 430.03845
Set-CellValue $ws 3 12 1780
Set-CellValue $ws 3 13 -316.03845
Set-CellValue $ws 3 14 -2008

# Sheet BSM, row 80
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 80 8 753.7692
Set-CellValue $ws 80 9 300.33334
Set-CellValue $ws 80 10 889.8
Set-CellValue $ws 80 11 300.33334
Set-CellValue $ws 80 12 889.8
Set-CellValue $ws 80 13 697.66666
Set-CellValue $ws 80 14 -2885.8

# Sheet BSM, row 83
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 83 8 753.7692
Set-CellValue $ws 83 9 300.33334
Set-CellValue $ws 83 10 889.8
Set-CellValue $ws 83 11 1501.6667
Set-CellValue $ws 83 12 4449
Set-CellValue $ws 83 13 3490.3333
Set-CellValue $ws 83 14 -14433

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 94 8 8064810.5
Set-CellValue $ws 94 9 8620985
Set-CellValue $ws 94 10 279
Set-CellValue $ws 94 11 8620985
Set-CellValue $ws 94 12 279
Set-CellValue $ws 94 13 -8620534
Set-CellValue $ws 94 14 -1181

# Sheet BSM, row 117
$ws = $wb.Worksheets.Item("BSM")
Set-CellValue $ws 117 8 30471.334
Set-CellValue $ws 117 10 30471.334
Set-CellValue $ws 117 12 30471.334
Set-CellValue $ws 117 14 -39649.334

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 16 8 55556772
Set-CellValue $ws 16 9 83334410
Set-CellValue $ws 16 10 1505.5
Set-CellValue $ws 16 11 83334410
Set-CellValue $ws 16 12 1505.5
Set-CellValue $ws 16 13 -83334123
Set-CellValue $ws 16 14 -2079.5

# Sheet CRP, row 22
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 22 8 47038.934
Set-CellValue $ws 22 9 252.75
Set-CellValue $ws 22 10 64052.09
Set-CellValue $ws 22 11 252.75
Set-CellValue $ws 22 12 64052.09
Set-CellValue $ws 22 13 97.25
Set-CellValue $ws 22 14 -64752.09

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 31 8 1373.1282
Set-CellValue $ws 31 9 1373.1282
Set-CellValue $ws 31 11 1373.1282
Set-CellValue $ws 31 13 -1078.1282

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 34 8 1373.1282
Set-CellValue $ws 34 9 1373.1282
Set-CellValue $ws 34 11 1373.1282
Set-CellValue $ws 34 13 -1171.1282

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 113 8 55556772
Set-CellValue $ws 113 9 83334410
Set-CellValue $ws 113 10 1505.5
Set-CellValue $ws 113 11 83334410
Set-CellValue $ws 113 12 1505.5
Set-CellValue $ws 113 13 -83332240
Set-CellValue $ws 113 14 -5845.5

# Sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
Set-CellValue $ws 122 8 1449.75
Set-CellValue $ws 122 9 1466.3334
Set-CellValue $ws 122 10 1400
Set-CellValue $ws 122 11 4399.0002
Set-CellValue $ws 122 12 4200
Set-CellValue $ws 122 13 -1949.0002
Set-CellValue $ws 122 14 -9100

# Sheet CUL, row 38
$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws 38 8 89.875
Set-CellValue $ws 38 9 62.545456
Set-CellValue $ws 38 10 150
Set-CellValue $ws 38 11 187.636368
Set-CellValue $ws 38 12 450
Set-CellValue $ws 38 13 159.363632
Set-CellValue $ws 38 14 -1144

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
Set-CellValue $ws 131 8 29415436
Set-CellValue $ws 131 9 76923540
Set-CellValue $ws 131 10 5660.5713
Set-CellValue $ws 131 11 230770620
Set-CellValue $ws 131 12 16981.7139
Set-CellValue $ws 131 13 -230765580
Set-CellValue $ws 131 14 -27061.7139

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
Set-CellValue $ws 126 8 1902.4736
Set-CellValue $ws 126 10 2089.6
Set-CellValue $ws 126 12 6268.799999999999
Set-CellValue $ws 126 14 -11208.8

# Sheet LTW, row 40
$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws 40 8 4133.8
Set-CellValue $ws 40 9 1791.4546
Set-CellValue $ws 40 11 1791.4546
Set-CellValue $ws 40 13 -1655.4546

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws 122 8 31251988
Set-CellValue $ws 122 9 41668332
Set-CellValue $ws 122 10 2950
Set-CellValue $ws 122 11 125004996
Set-CellValue $ws 122 12 8850
Set-CellValue $ws 122 13 -125002546
Set-CellValue $ws 122 14 -13750

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
Set-CellValue $ws 132 8 2714.72
Set-CellValue $ws 132 9 2525.5334
Set-CellValue $ws 132 10 2998.5
Set-CellValue $ws 132 11 7576.600199999999
Set-CellValue $ws 132 12 8995.5
Set-CellValue $ws 132 13 -5046.600199999999
Set-CellValue $ws 132 14 -14055.5

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws 122 8 9260719
Set-CellValue $ws 122 9 9616823
Set-CellValue $ws 122 11 28850469
Set-CellValue $ws 122 13 -28848019

# Sheet WVR, row 137
$ws = $wb.Worksheets.Item("WVR")
Set-CellValue $ws 137 8 33053.332
Set-CellValue $ws 137 10 33053.332
Set-CellValue $ws 137 12 33053.332
Set-CellValue $ws 137 14 -43253.332
